# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the upstream data source (gh-pages
# regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 3107
$ws1.Range("F3").Value  = 509
$ws1.Range("F5").Value  = 69
$ws1.Range("F6").Value  = 17
$ws1.Range("F8").Value  = 23
$ws1.Range("F9").Value  = 1093
$ws1.Range("F10").Value = 15252
$ws1.Range("F11").Value = 212
$ws1.Range("F13").Value = 1017
$ws1.Range("F14").Value = 6038
$ws1.Range("F16").Value = 96
$ws1.Range("F18").Value = 98
$ws1.Range("F26").Value = 4970
$ws1.Range("F27").Value = 127
$ws1.Range("F28").Value = 10898
$ws1.Range("F31").Value = 102
$ws1.Range("F32").Value = 145
$ws1.Range("F33").Value = 3777

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 3107
$ws4.Range("F4").Value  = 509
$ws4.Range("F6").Value  = 69
$ws4.Range("F7").Value  = 17
$ws4.Range("F9").Value  = 23
$ws4.Range("F10").Value = 1093
$ws4.Range("F11").Value = 15252
$ws4.Range("F12").Value = 212
$ws4.Range("F14").Value = 1017
$ws4.Range("F15").Value = 6038
$ws4.Range("F17").Value = 96
$ws4.Range("F19").Value = 98
$ws4.Range("F27").Value = 4970
$ws4.Range("F28").Value = 127
$ws4.Range("F30").Value = 10898
$ws4.Range("F33").Value = 102
$ws4.Range("F34").Value = 145
$ws4.Range("F35").Value = 3777

$wb.Save()
